$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3.627806666666667
$ws.Cells.Item(2, 8).Value = 10.88342
$ws.Cells.Item(2, 9).Value = 0.08036675778141429
$ws.Cells.Item(2, 10).Value = 0.08036675778141429
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.731629
$ws.Cells.Item(2, 14).Value = 8.194887
$ws.Cells.Item(2, 15).Value = 0.5547800938501829
$ws.Cells.Item(2, 16).Value = 0.554780093850183
$ws.Cells.Item(2, 17).Value = 9.90982189706
$ws.Cells.Item(2, 18).Value = 89.18839707354
$ws.Cells.Item(2, 19).Value = 0.04458587742440794
$ws.Cells.Item(2, 20).Value = 0.04458587742440795
$ws.Cells.Item(3, 7).Value = 3.627806666666667
$ws.Cells.Item(3, 8).Value = 10.88342
$ws.Cells.Item(3, 9).Value = 0.08036675778141429
$ws.Cells.Item(3, 10).Value = 0.08036675778141429
$ws.Cells.Item(3, 13).Value = 0.06813733333333333
$ws.Cells.Item(3, 15).Value = 0.01383834927121065
$ws.Cells.Item(3, 16).Value = 0.01383834927121065
$ws.Cells.Item(3, 17).Value = 0.2471890721155555
$ws.Cells.Item(3, 18).Value = 2.22470164904
$ws.Cells.Item(3, 19).Value = 0.001112143263973997
$ws.Cells.Item(3, 20).Value = 0.001112143263973997
$ws.Cells.Item(4, 7).Value = 3.627806666666667
$ws.Cells.Item(4, 8).Value = 10.88342
$ws.Cells.Item(4, 9).Value = 0.08036675778141429
$ws.Cells.Item(4, 10).Value = 0.08036675778141429
$ws.Cells.Item(4, 13).Value = 2.124038666666666
$ws.Cells.Item(4, 14).Value = 6.372115999999999
$ws.Cells.Item(4, 15).Value = 0.4313815568786064
$ws.Cells.Item(4, 16).Value = 0.4313815568786064
$ws.Cells.Item(4, 17).Value = 7.705601635191111
$ws.Cells.Item(4, 18).Value = 69.35041471672
$ws.Cells.Item(4, 19).Value = 0.03466873709303235
$ws.Cells.Item(4, 20).Value = 0.03466873709303235
$ws.Cells.Item(5, 9).Value = 0.6506403335968259
$ws.Cells.Item(5, 10).Value = 0.6506403335968259
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.731629
$ws.Cells.Item(5, 14).Value = 8.194887
$ws.Cells.Item(5, 15).Value = 0.5547800938501829
$ws.Cells.Item(5, 16).Value = 0.554780093850183
$ws.Cells.Item(5, 17).Value = 80.22881603019398
$ws.Cells.Item(5, 18).Value = 722.0593442717459
$ws.Cells.Item(5, 19).Value = 0.3609623053355614
$ws.Cells.Item(5, 20).Value = 0.3609623053355615
$ws.Cells.Item(6, 9).Value = 0.6506403335968259
$ws.Cells.Item(6, 10).Value = 0.6506403335968259
$ws.Cells.Item(6, 13).Value = 0.06813733333333333
$ws.Cells.Item(6, 15).Value = 0.01383834927121065
$ws.Cells.Item(6, 16).Value = 0.01383834927121065
$ws.Cells.Item(6, 19).Value = 0.009003788186249887
$ws.Cells.Item(6, 20).Value = 0.009003788186249889
$ws.Cells.Item(7, 9).Value = 0.6506403335968259
$ws.Cells.Item(7, 10).Value = 0.6506403335968259
$ws.Cells.Item(7, 13).Value = 2.124038666666666
$ws.Cells.Item(7, 14).Value = 6.372115999999999
$ws.Cells.Item(7, 15).Value = 0.4313815568786064
$ws.Cells.Item(7, 16).Value = 0.4313815568786064
$ws.Cells.Item(7, 17).Value = 62.38369391634755
$ws.Cells.Item(7, 18).Value = 561.4532452471279
$ws.Cells.Item(7, 19).Value = 0.2806742400750146
$ws.Cells.Item(7, 20).Value = 0.2806742400750146
$ws.Cells.Item(8, 7).Value = 12.14251133333333
$ws.Cells.Item(8, 8).Value = 36.427534
$ws.Cells.Item(8, 9).Value = 0.2689929086217598
$ws.Cells.Item(8, 10).Value = 0.2689929086217598
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.731629
$ws.Cells.Item(8, 14).Value = 8.194887
$ws.Cells.Item(8, 15).Value = 0.5547800938501829
$ws.Cells.Item(8, 16).Value = 0.554780093850183
$ws.Cells.Item(8, 17).Value = 33.168836090962
$ws.Cells.Item(8, 18).Value = 298.519524818658
$ws.Cells.Item(8, 19).Value = 0.1492319110902136
$ws.Cells.Item(8, 20).Value = 0.1492319110902136
$ws.Cells.Item(9, 7).Value = 12.14251133333333
$ws.Cells.Item(9, 8).Value = 36.427534
$ws.Cells.Item(9, 9).Value = 0.2689929086217598
$ws.Cells.Item(9, 10).Value = 0.2689929086217598
$ws.Cells.Item(9, 13).Value = 0.06813733333333333
$ws.Cells.Item(9, 15).Value = 0.01383834927121065
$ws.Cells.Item(9, 16).Value = 0.01383834927121065
$ws.Cells.Item(9, 17).Value = 0.827358342223111
$ws.Cells.Item(9, 18).Value = 7.446225080007999
$ws.Cells.Item(9, 19).Value = 0.003722417820986762
$ws.Cells.Item(9, 20).Value = 0.003722417820986763
$ws.Cells.Item(10, 7).Value = 12.14251133333333
$ws.Cells.Item(10, 8).Value = 36.427534
$ws.Cells.Item(10, 9).Value = 0.2689929086217598
$ws.Cells.Item(10, 10).Value = 0.2689929086217598
$ws.Cells.Item(10, 13).Value = 2.124038666666666
$ws.Cells.Item(10, 14).Value = 6.372115999999999
$ws.Cells.Item(10, 15).Value = 0.4313815568786064
$ws.Cells.Item(10, 16).Value = 0.4313815568786064
$ws.Cells.Item(10, 17).Value = 25.79116358243822
$ws.Cells.Item(10, 18).Value = 232.120472241944
$ws.Cells.Item(10, 19).Value = 0.1160385797105595
$ws.Cells.Item(10, 20).Value = 0.1160385797105595
